# Auto-generated edit script applying the Daily Update (10 keyword rotation) diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row=3; Cells=@{ E="188670" } },
    @{ Row=12; Cells=@{ B="벤큐 GW2780 아이케어"; C="https://search.shopping.naver.com/gate.nhn?id=11846746201"; D="https://shopping-phinf.pstatic.net/main_1184674/11846746201.20211118104344.jpg"; E="209000" } },
    @{ Row=13; Cells=@{ B="LG전자 27TN600S"; C="https://search.shopping.naver.com/gate.nhn?id=21745903830"; D="https://shopping-phinf.pstatic.net/main_2174590/21745903830.20200824150453.jpg"; E="341800"; I="LG전자"; J="LG전자" } },
    @{ Row=14; Cells=@{ B="벤큐 XL2411K"; C="https://search.shopping.naver.com/gate.nhn?id=24196764522"; D="https://shopping-phinf.pstatic.net/main_2419676/24196764522.20211117182333.jpg"; E="279000"; I="벤큐"; J="벤큐" } },
    @{ Row=17; Cells=@{ B="삼성전자 스마트모니터 M7 S43AM700"; C="https://search.shopping.naver.com/gate.nhn?id=27215825524"; D="https://shopping-phinf.pstatic.net/main_2721582/27215825524.20210520173025.jpg"; E="649000" } },
    @{ Row=18; Cells=@{ B="ASUS VA24DQLB"; C="https://search.shopping.naver.com/gate.nhn?id=23257141490"; D="https://shopping-phinf.pstatic.net/main_2325714/23257141490.20200720175614.jpg"; E="219000"; I="ASUS"; J="ASUS" } },
    @{ Row=20; Cells=@{ B="삼성전자 스마트모니터 M7 S32AM700"; C="https://search.shopping.naver.com/gate.nhn?id=25524333522"; D="https://shopping-phinf.pstatic.net/main_2552433/25524333522.20210203133321.jpg"; E="499000"; I="스마트모니터"; J="삼성전자" } },
    @{ Row=36; Cells=@{ B="삼성전자 오디세이 G5 C32G54T"; C="https://search.shopping.naver.com/gate.nhn?id=23896004523"; D="https://shopping-phinf.pstatic.net/main_2389600/23896004523.20210203132926.jpg"; E="420000"; I="오디세이"; J="삼성전자" } },
    @{ Row=37; Cells=@{ B="한성컴퓨터 TFG32Q07P 75"; C="https://search.shopping.naver.com/gate.nhn?id=28655748554"; D="https://shopping-phinf.pstatic.net/main_2865574/28655748554.20210831152013.jpg"; E="259000"; I="한성컴퓨터"; J="한성컴퓨터" } },
    @{ Row=61; Cells=@{ E="237000" } },
    @{ Row=64; Cells=@{ B="삼성전자 오디세이 G9 C49G95T"; C="https://search.shopping.naver.com/gate.nhn?id=23255316490"; D="https://shopping-phinf.pstatic.net/main_2325531/23255316490.20210203133605.jpg"; E="1690000"; I="오디세이"; J="삼성전자" } },
    @{ Row=65; Cells=@{ B="벤큐 XL2546K"; C="https://search.shopping.naver.com/gate.nhn?id=24235203522"; D="https://shopping-phinf.pstatic.net/main_2423520/24235203522.20211117182410.jpg"; E="649000"; I="벤큐"; J="벤큐" } },
    @{ Row=67; Cells=@{ B="알파스캔 AOC 27B2 보더리스 75 시력보호"; C="https://search.shopping.naver.com/gate.nhn?id=21720504796"; D="https://shopping-phinf.pstatic.net/main_2172050/21720504796.20210310171806.jpg"; E="219000"; I="알파스캔"; J="알파스캔" } },
    @{ Row=68; Cells=@{ B="한성컴퓨터 TFG39Q14V 144"; C="https://search.shopping.naver.com/gate.nhn?id=26826361522"; D="https://shopping-phinf.pstatic.net/main_2682636/26826361522.20210419161946.jpg"; E="499000"; I="한성컴퓨터"; J="한성컴퓨터" } },
    @{ Row=69; Cells=@{ B="벤큐 ZOWIE XL2731"; C="https://search.shopping.naver.com/gate.nhn?id=22435628535"; D="https://shopping-phinf.pstatic.net/main_2243562/22435628535.20211126161127.jpg"; E="419000"; I="벤큐"; J="벤큐" } },
    @{ Row=73; Cells=@{ E="567000" } },
    @{ Row=91; Cells=@{ B="LG전자 울트라와이드 29WP500"; C="https://search.shopping.naver.com/gate.nhn?id=26886077522"; D="https://shopping-phinf.pstatic.net/main_2688607/26886077522.20210524134552.jpg"; E="259000"; I="울트라와이드"; J="LG전자" } },
    @{ Row=92; Cells=@{ B="한성컴퓨터 TFG27Q14P 144"; C="https://search.shopping.naver.com/gate.nhn?id=27327723522"; D="https://shopping-phinf.pstatic.net/main_2732772/27327723522.20210527095004.jpg"; E="379000"; I="한성컴퓨터"; J="한성컴퓨터" } },
    @{ Row=93; Cells=@{ B="삼성전자 삼성 U32R590"; C="https://search.shopping.naver.com/gate.nhn?id=17650306747"; D="https://shopping-phinf.pstatic.net/main_1765030/17650306747.20210203134432.jpg"; E="399000"; I="삼성"; J="삼성전자" } },
    @{ Row=94; Cells=@{ B="삼성전자 삼성 C27F391"; C="https://search.shopping.naver.com/gate.nhn?id=9681100715"; D="https://shopping-phinf.pstatic.net/main_9681100/9681100715.20200915114554.jpg"; E="238000"; I="삼성"; J="삼성전자" } },
    @{ Row=95; Cells=@{ B="DELL 울트라샤프 U2720Q"; C="https://search.shopping.naver.com/gate.nhn?id=21752731630"; D="https://shopping-phinf.pstatic.net/main_2175273/21752731630.20200327122054.jpg"; E="778990"; I="울트라샤프"; J="DELL" } },
    @{ Row=100; Cells=@{ B="삼성전자 삼성 C27F390"; C="https://search.shopping.naver.com/gate.nhn?id=9489557554"; D="https://shopping-phinf.pstatic.net/main_9489557/9489557554.20210203132811.jpg"; E="238000"; I="삼성"; J="삼성전자" } },
    @{ Row=101; Cells=@{ B="벤큐 모비우스 EX2710S"; C="https://search.shopping.naver.com/gate.nhn?id=27862189523"; D="https://shopping-phinf.pstatic.net/main_2786218/27862189523.20211116103224.jpg"; E="369000"; I="벤큐"; J="벤큐" } }
)

foreach ($item in $changes) {
    $r = $item.Row
    foreach ($col in $item.Cells.Keys) {
        $addr = "$col$r"
        $val = $item.Cells[$col]
        if ($col -eq "E") {
            # Price column is stored as text in the source data; force text format
            # so the numeric-looking string is not auto-converted to a Number.
            $ws.Range($addr).NumberFormat = "@"
        }
        $ws.Range($addr).Value = $val
    }
}
